$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the three "early warning" notes to "withdrew, early warning"
# (the student in D23, "SAS accomodations", is left unchanged)
$ws.Range("D5").Value = "withdrew, early warning"
$ws.Range("D9").Value = "withdrew, early warning"
$ws.Range("D30").Value = "withdrew, early warning"

# Scroll the view down and select the whole used range, as recorded
# by Excel when the workbook was last saved.
$ws.Application.ActiveWindow.ScrollRow = 20
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Range("A1:D30").Select()
